$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 5.403619666666667
$ws.Range("H2").Value = 16.210859
$ws.Range("I2").Value = 0.08303620947547638
$ws.Range("J2").Value = 0.08303620947547637
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.1575256666666667
$ws.Range("N2").Value = 0.472577
$ws.Range("O2").Value = 0.6985926944284299
$ws.Range("P2").Value = 0.69859269442843
$ws.Range("Q2").Value = 0.8512087904047778
$ws.Range("R2").Value = 7.660879113643
$ws.Range("S2").Value = 0.05800848931259656
$ws.Range("T2").Value = 0.05800848931259656

$ws.Range("G3").Value = 5.403619666666667
$ws.Range("H3").Value = 16.210859
$ws.Range("I3").Value = 0.08303620947547638
$ws.Range("J3").Value = 0.08303620947547637
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.06796433333333333
$ws.Range("N3").Value = 0.203893
$ws.Range("O3").Value = 0.30140730557157
$ws.Range("P3").Value = 0.30140730557157
$ws.Range("Q3").Value = 0.3672534082318889
$ws.Range("R3").Value = 3.305280674087
$ws.Range("S3").Value = 0.02502772016287981
$ws.Range("T3").Value = 0.0250277201628798

$ws.Range("G4").Value = 50.79415266666667
$ws.Range("H4").Value = 152.382458
$ws.Range("I4").Value = 0.7805423329433673
$ws.Range("J4").Value = 0.7805423329433673
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.1575256666666667
$ws.Range("N4").Value = 0.472577
$ws.Range("O4").Value = 0.6985926944284299
$ws.Range("P4").Value = 0.69859269442843
$ws.Range("Q4").Value = 8.001382761585113
$ws.Range("R4").Value = 72.01244485426601
$ws.Range("S4").Value = 0.5452811714863596
$ws.Range("T4").Value = 0.5452811714863596

$ws.Range("G5").Value = 50.79415266666667
$ws.Range("H5").Value = 152.382458
$ws.Range("I5").Value = 0.7805423329433673
$ws.Range("J5").Value = 0.7805423329433673
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.06796433333333333
$ws.Range("N5").Value = 0.203893
$ws.Range("O5").Value = 0.30140730557157
$ws.Range("P5").Value = 0.30140730557157
$ws.Range("Q5").Value = 3.452190723221556
$ws.Range("R5").Value = 31.069716508994
$ws.Range("S5").Value = 0.2352611614570077
$ws.Range("T5").Value = 0.2352611614570077

$ws.Range("G6").Value = 8.877689333333334
$ws.Range("H6").Value = 26.633068
$ws.Range("I6").Value = 0.1364214575811564
$ws.Range("J6").Value = 0.1364214575811563
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.1575256666666667
$ws.Range("N6").Value = 0.472577
$ws.Range("O6").Value = 0.6985926944284299
$ws.Range("P6").Value = 0.69859269442843
$ws.Range("Q6").Value = 1.398463930692889
$ws.Range("R6").Value = 12.586175376236
$ws.Range("S6").Value = 0.09530303362947377
$ws.Range("T6").Value = 0.09530303362947377

$ws.Range("G7").Value = 8.877689333333334
$ws.Range("H7").Value = 26.633068
$ws.Range("I7").Value = 0.1364214575811564
$ws.Range("J7").Value = 0.1364214575811563
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.06796433333333333
$ws.Range("N7").Value = 0.203893
$ws.Range("O7").Value = 0.30140730557157
$ws.Range("P7").Value = 0.30140730557157
$ws.Range("Q7").Value = 0.6033662370804446
$ws.Range("R7").Value = 5.430296133724
$ws.Range("S7").Value = 0.04111842395168257
$ws.Range("T7").Value = 0.04111842395168257
